$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing E-column values (daily refresh of outbreak counts) ---
$ws.Range("E11").Value = 419
$ws.Range("E28").Value = 4781
$ws.Range("E48").Value = 19270
$ws.Range("E62").Value = 7
$ws.Range("E69").Value = 28517
$ws.Range("E71").Value = 2076
$ws.Range("E90").Value = 27654
$ws.Range("E111").Value = 19060
$ws.Range("E114").Value = 3242
$ws.Range("E124").Value = 192
$ws.Range("E132").Value = 11766
$ws.Range("E135").Value = 2246
$ws.Range("E136").Value = 532
$ws.Range("E151").Value = 8224
$ws.Range("E155").Value = 542
$ws.Range("E168").Value = 4891
$ws.Range("E171").Value = 706
$ws.Range("E184").Value = 4046
$ws.Range("E201").Value = 2970
$ws.Range("E205").Value = 307
$ws.Range("E223").Value = 341
$ws.Range("E245").Value = 20
$ws.Range("E357").Value = 199
$ws.Range("E366").Value = 1832
$ws.Range("E367").Value = 27
$ws.Range("E386").Value = 2283
$ws.Range("E387").Value = 38
$ws.Range("E388").Value = 663
$ws.Range("E406").Value = 2949
$ws.Range("E424").Value = 4075
$ws.Range("E426").Value = 919
$ws.Range("E443").Value = 273
$ws.Range("E462").Value = 7085
$ws.Range("E464").Value = 1245
$ws.Range("E479").Value = 6555
$ws.Range("E481").Value = 1060
$ws.Range("E500").Value = 1174
$ws.Range("E518").Value = 6990
$ws.Range("E520").Value = 1280
$ws.Range("E538").Value = 382
$ws.Range("E539").Value = 8693
$ws.Range("E541").Value = 1445
$ws.Range("E558").Value = 525
$ws.Range("E559").Value = 9398
$ws.Range("E574").Value = 299
$ws.Range("E580").Value = 11382
$ws.Range("E582").Value = 1891
$ws.Range("E583").Value = 265
$ws.Range("E601").Value = 19469
$ws.Range("E603").Value = 2813
$ws.Range("E622").Value = 32903
$ws.Range("E624").Value = 3803
$ws.Range("E632").Value = 734
$ws.Range("E643").Value = 61721
$ws.Range("E645").Value = 5441
$ws.Range("E647").Value = 280
$ws.Range("E651").Value = 56
$ws.Range("E653").Value = 1112
$ws.Range("E664").Value = 94722
$ws.Range("E666").Value = 7114
$ws.Range("E667").Value = 2466
$ws.Range("E668").Value = 358
$ws.Range("E685").Value = 108310
$ws.Range("E687").Value = 8104
$ws.Range("E688").Value = 2992
$ws.Range("E693").Value = 98
$ws.Range("E696").Value = 398
$ws.Range("E706").Value = 108310
$ws.Range("E709").Value = 4517
$ws.Range("E717").Value = 536
$ws.Range("E727").Value = 106254
$ws.Range("E729").Value = 7759
$ws.Range("E730").Value = 6131
$ws.Range("E734").Value = 818
$ws.Range("E737").Value = 1064
$ws.Range("E738").Value = 747
$ws.Range("E739").Value = 332
$ws.Range("E740").Value = 328
$ws.Range("E741").Value = 502
$ws.Range("E748").Value = 100754
$ws.Range("E750").Value = 7105
$ws.Range("E751").Value = 6784
$ws.Range("E759").Value = 673
$ws.Range("E761").Value = 422
$ws.Range("E769").Value = 104644
$ws.Range("E771").Value = 6618
$ws.Range("E772").Value = 7758
$ws.Range("E773").Value = 432
$ws.Range("E776").Value = 1383
$ws.Range("E779").Value = 1232
$ws.Range("E780").Value = 940
$ws.Range("E781").Value = 408
$ws.Range("E782").Value = 587
$ws.Range("E783").Value = 424
$ws.Range("E788").Value = 548
$ws.Range("E790").Value = 126429
$ws.Range("E792").Value = 7549
$ws.Range("E793").Value = 10192
$ws.Range("E794").Value = 511
$ws.Range("E797").Value = 2114
$ws.Range("E800").Value = 1537
$ws.Range("E802").Value = 806
$ws.Range("E803").Value = 891
$ws.Range("E809").Value = 595
$ws.Range("E810").Value = 3272
$ws.Range("E811").Value = 143489
$ws.Range("E813").Value = 7502
$ws.Range("E814").Value = 11894
$ws.Range("E817").Value = 117
$ws.Range("E818").Value = 2140
$ws.Range("E819").Value = 149
$ws.Range("E821").Value = 1442
$ws.Range("E824").Value = 792
$ws.Range("E825").Value = 667
$ws.Range("E831").Value = 115525
$ws.Range("E833").Value = 5237
$ws.Range("E834").Value = 10264
$ws.Range("E841").Value = 733
$ws.Range("E844").Value = 529
$ws.Range("E849").Value = 389
$ws.Range("E850").Value = 2480
$ws.Range("E851").Value = 99541
$ws.Range("E853").Value = 6415
$ws.Range("E854").Value = 10792
$ws.Range("E858").Value = 1369
$ws.Range("E869").Value = 380
$ws.Range("E870").Value = 1890
$ws.Range("E871").Value = 118061
$ws.Range("E873").Value = 8430
$ws.Range("E874").Value = 11790
$ws.Range("E877").Value = 72
$ws.Range("E878").Value = 1725
$ws.Range("E881").Value = 661
$ws.Range("E890").Value = 94085
$ws.Range("E892").Value = 6607
$ws.Range("E893").Value = 10972
$ws.Range("E896").Value = 107
$ws.Range("E897").Value = 1906
$ws.Range("E900").Value = 814
$ws.Range("E903").Value = 301
$ws.Range("E908").Value = 330
$ws.Range("E909").Value = 1975
$ws.Range("E910").Value = 73520
$ws.Range("E912").Value = 5523
$ws.Range("E913").Value = 9075
$ws.Range("E917").Value = 1855
$ws.Range("E922").Value = 284
$ws.Range("E923").Value = 303
$ws.Range("E929").Value = 379
$ws.Range("E930").Value = 1697
$ws.Range("E931").Value = 59885
$ws.Range("E933").Value = 5086
$ws.Range("E934").Value = 6404
$ws.Range("E937").Value = 63
$ws.Range("E938").Value = 1596
$ws.Range("E941").Value = 1402
$ws.Range("E945").Value = 341
$ws.Range("E950").Value = 253
$ws.Range("E952").Value = 49860
$ws.Range("E953").Value = 223
$ws.Range("E954").Value = 4676
$ws.Range("E955").Value = 3963
$ws.Range("E959").Value = 1346
$ws.Range("E962").Value = 1616
$ws.Range("E973").Value = 38574
$ws.Range("E975").Value = 4515
$ws.Range("E976").Value = 2456
$ws.Range("E979").Value = 41
$ws.Range("E980").Value = 1055
$ws.Range("E982").Value = 72
$ws.Range("E983").Value = 1651
$ws.Range("E993").Value = 808
$ws.Range("E994").Value = 39633
$ws.Range("E996").Value = 4984
$ws.Range("E997").Value = 1792
$ws.Range("E1003").Value = 117
$ws.Range("E1008").Value = 523
$ws.Range("E1014").Value = 1018
$ws.Range("E1015").Value = 42705
$ws.Range("E1017").Value = 5818
$ws.Range("E1018").Value = 1293
$ws.Range("E1019").Value = 204
$ws.Range("E1022").Value = 763
$ws.Range("E1024").Value = 87
$ws.Range("E1026").Value = 136
$ws.Range("E1034").Value = 350
$ws.Range("E1035").Value = 1099
$ws.Range("E1036").Value = 45598
$ws.Range("E1037").Value = 210
$ws.Range("E1038").Value = 5925
$ws.Range("E1044").Value = 172
$ws.Range("E1045").Value = 56
$ws.Range("E1046").Value = 1931
$ws.Range("E1047").Value = 301
$ws.Range("E1048").Value = 131
$ws.Range("E1049").Value = 185
$ws.Range("E1050").Value = 1055
$ws.Range("E1054").Value = 359
$ws.Range("E1055").Value = 964
$ws.Range("E1056").Value = 56706
$ws.Range("E1057").Value = 271
$ws.Range("E1058").Value = 6953
$ws.Range("E1059").Value = 709
$ws.Range("E1060").Value = 218
$ws.Range("E1062").Value = 57
$ws.Range("E1063").Value = 385
$ws.Range("E1064").Value = 100
$ws.Range("E1066").Value = 1981
$ws.Range("E1067").Value = 478
$ws.Range("E1068").Value = 172
$ws.Range("E1069").Value = 80
$ws.Range("E1070").Value = 1710
$ws.Range("E1071").Value = 131
$ws.Range("E1073").Value = 297
$ws.Range("E1074").Value = 1061
$ws.Range("E1075").Value = 76490
$ws.Range("E1077").Value = 7643
$ws.Range("E1078").Value = 735
$ws.Range("E1079").Value = 238
$ws.Range("E1080").Value = 173
$ws.Range("E1081").Value = 41
$ws.Range("E1082").Value = 397
$ws.Range("E1083").Value = 153
$ws.Range("E1084").Value = 42
$ws.Range("E1085").Value = 2123
$ws.Range("E1086").Value = 412
$ws.Range("E1087").Value = 175
$ws.Range("E1088").Value = 155
$ws.Range("E1089").Value = 1775
$ws.Range("E1090").Value = 145
$ws.Range("E1094").Value = 410
$ws.Range("E1095").Value = 1354
$ws.Range("E1096").Value = 100364
$ws.Range("E1097").Value = 296
$ws.Range("E1098").Value = 7966
$ws.Range("E1099").Value = 630
$ws.Range("E1100").Value = 302
$ws.Range("E1101").Value = 155
$ws.Range("E1102").Value = 55
$ws.Range("E1103").Value = 274
$ws.Range("E1104").Value = 82
$ws.Range("E1105").Value = 48
$ws.Range("E1106").Value = 2005
$ws.Range("E1107").Value = 481
$ws.Range("E1108").Value = 208
$ws.Range("E1109").Value = 71
$ws.Range("E1110").Value = 1702
$ws.Range("E1111").Value = 177
$ws.Range("E1113").Value = 24
$ws.Range("E1114").Value = 8
$ws.Range("E1115").Value = 372
$ws.Range("E1116").Value = 1405

# --- Append new rows for Meldewoche 13 ---
$ws.Range("A1117").Value = 2021
$ws.Range("B1117").Value = 13
$ws.Range("C1117").Value = "Nicht in Ausbruch erfasst"
$ws.Range("D1117").Value = "Not documented in an outbreak"
$ws.Range("E1117").Value = 99224
$ws.Range("A1118").Value = 2021
$ws.Range("B1118").Value = 13
$ws.Range("C1118").Value = "Wohnstätten"
$ws.Range("D1118").Value = "Residences"
$ws.Range("E1118").Value = 184
$ws.Range("A1119").Value = 2021
$ws.Range("B1119").Value = 13
$ws.Range("C1119").Value = "Privater Haushalt"
$ws.Range("D1119").Value = "Private household"
$ws.Range("E1119").Value = 5092
$ws.Range("A1120").Value = 2021
$ws.Range("B1120").Value = 13
$ws.Range("C1120").Value = "Alten-/Pflegeheim"
$ws.Range("D1120").Value = "Retirement/nursing home"
$ws.Range("E1120").Value = 459
$ws.Range("A1121").Value = 2021
$ws.Range("B1121").Value = 13
$ws.Range("C1121").Value = "Flüchtlingsheim"
$ws.Range("D1121").Value = "Refugee accomodation"
$ws.Range("E1121").Value = 224
$ws.Range("A1122").Value = 2021
$ws.Range("B1122").Value = 13
$ws.Range("C1122").Value = "Wohnheim"
$ws.Range("D1122").Value = "Residential home"
$ws.Range("E1122").Value = 49
$ws.Range("A1123").Value = 2021
$ws.Range("B1123").Value = 13
$ws.Range("C1123").Value = "Med. Behandlungseinrichtung"
$ws.Range("D1123").Value = "Health care centre"
$ws.Range("E1123").Value = 27
$ws.Range("A1124").Value = 2021
$ws.Range("B1124").Value = 13
$ws.Range("C1124").Value = "Krankenhaus"
$ws.Range("D1124").Value = "Hospital"
$ws.Range("E1124").Value = 145
$ws.Range("A1125").Value = 2021
$ws.Range("B1125").Value = 13
$ws.Range("C1125").Value = "Praxis"
$ws.Range("D1125").Value = "Medical practice"
$ws.Range("E1125").Value = 38
$ws.Range("A1126").Value = 2021
$ws.Range("B1126").Value = 13
$ws.Range("C1126").Value = "Reha-Einrichtung"
$ws.Range("D1126").Value = "Medical rehabilitation"
$ws.Range("E1126").Value = 15
$ws.Range("A1127").Value = 2021
$ws.Range("B1127").Value = 13
$ws.Range("C1127").Value = "Arbeitsplatz"
$ws.Range("D1127").Value = "Work place"
$ws.Range("E1127").Value = 1036
$ws.Range("A1128").Value = 2021
$ws.Range("B1128").Value = 13
$ws.Range("C1128").Value = "Ausbildungsstätte"
$ws.Range("D1128").Value = "Educational institution"
$ws.Range("E1128").Value = 216
$ws.Range("A1129").Value = 2021
$ws.Range("B1129").Value = 13
$ws.Range("C1129").Value = "Betreuungseinrichtung"
$ws.Range("D1129").Value = "Care facility"
$ws.Range("E1129").Value = 92
$ws.Range("A1130").Value = 2021
$ws.Range("B1130").Value = 13
$ws.Range("C1130").Value = "Seniorentagesstätte"
$ws.Range("D1130").Value = "Day-care centre for the elderly"
$ws.Range("E1130").Value = 26
$ws.Range("A1131").Value = 2021
$ws.Range("B1131").Value = 13
$ws.Range("C1131").Value = "Kindergarten, Hort"
$ws.Range("D1131").Value = "Kindergarten, after-school child care"
$ws.Range("E1131").Value = 684
$ws.Range("A1132").Value = 2021
$ws.Range("B1132").Value = 13
$ws.Range("C1132").Value = "Freizeit"
$ws.Range("D1132").Value = "Leisure"
$ws.Range("E1132").Value = 87
$ws.Range("A1133").Value = 2021
$ws.Range("B1133").Value = 13
$ws.Range("C1133").Value = "Speisestätte"
$ws.Range("D1133").Value = "Dining venue"
$ws.Range("E1133").Value = 6
$ws.Range("A1134").Value = 2021
$ws.Range("B1134").Value = 13
$ws.Range("C1134").Value = "Übernachtung"
$ws.Range("D1134").Value = "Overnight stay"
$ws.Range("E1134").Value = 3
$ws.Range("A1135").Value = 2021
$ws.Range("B1135").Value = 13
$ws.Range("C1135").Value = "Verkehrsmittel"
$ws.Range("D1135").Value = "Public transport"
$ws.Range("E1135").Value = 8
$ws.Range("A1136").Value = 2021
$ws.Range("B1136").Value = 13
$ws.Range("C1136").Value = "weitere Settings"
$ws.Range("D1136").Value = "Other"
$ws.Range("E1136").Value = 207
$ws.Range("A1137").Value = 2021
$ws.Range("B1137").Value = 13
$ws.Range("C1137").Value = "unbekanntes Setting"
$ws.Range("D1137").Value = "Unknown"
$ws.Range("E1137").Value = 916
